$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the values: E10 and E11 change from "ACTUALIZADO" to "n/a"
$ws.Range("E10").Value = "n/a"
$ws.Range("E11").Value = "n/a"

# Update the sheet view: scroll so column D is the left-most visible column
# (topLeftCell="D1") and move the selection to E12.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("E12").Select()
